# Generate Report for Handback
# Update the "Correspond Handback Datetime" values for row 3 (the
# 8c9ec3f7-... file) on both the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-24 19:02:51"
$wsZhCn.Range("H3").Value = "2016-03-24 19:03:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-24 19:02:55"
$wsDeDe.Range("H3").Value = "2016-03-24 19:03:37"
